$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "SSB and Bmsy data"
$ws2 = $wb.Worksheets.Item(2)   # "summary SSB _Bmsy per spp"

# --- Active sheet / tab selection -----------------------------------------
# The commit moves the "active"/selected tab from the summary sheet back to
# the "SSB and Bmsy data" sheet.
$null = $ws1.Activate()

# --- Column width on "SSB and Bmsy data" -----------------------------------
# Column D width changes from 19.6640625 to 17.1640625 (character units).
$ws1.Columns.Item(4).ColumnWidth = 16.330729166666668

# --- New highlight fill on header row (row 8) ------------------------------
# C8, F8, H8, L8 and N8 get a new (4th) fill style - a light themed orange
# highlight - applied via Interior color.
$highlightCells = @("C8", "F8", "H8", "L8", "N8")
foreach ($addr in $highlightCells) {
    $ws1.Range($addr).Interior.Color = 14083579
}

# --- Selection / view state on "SSB and Bmsy data" --------------------------
# Selection moves from B3:M4 to K6 (inside the frozen bottom pane).
# (Do this last - selecting a range activates its parent sheet, and we need
# "SSB and Bmsy data" to stay the active/selected tab, matching the diff.)
$null = $ws1.Range("K6").Select()
